$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report generation timestamp in A1
$ws.Range("A1").Value = "CreatedAt: 2025-11-04T19:06:34"

# Update intertie LMP component values in columns W:Z (hours 21-24)
$ws.Range("W4").Value = 88.92
$ws.Range("X4").Value = 74.20999999999999
$ws.Range("Y4").Value = 91.89
$ws.Range("Z4").Value = 74.23999999999999
$ws.Range("W6").Value = -6.14
$ws.Range("X6").Value = -4.9
$ws.Range("Y6").Value = -5.88
$ws.Range("Z6").Value = -3.64
$ws.Range("W9").Value = 92.02
$ws.Range("X9").Value = 77.55
$ws.Range("Y9").Value = 97.09
$ws.Range("Z9").Value = 79.45999999999999
$ws.Range("W11").Value = -3.04
$ws.Range("X11").Value = -1.55
$ws.Range("Y11").Value = -0.68
$ws.Range("Z11").Value = 1.59
$ws.Range("W14").Value = 92.02
$ws.Range("X14").Value = 77.55
$ws.Range("Y14").Value = 97.09
$ws.Range("Z14").Value = 79.54000000000001
$ws.Range("W16").Value = -3.04
$ws.Range("X16").Value = -1.55
$ws.Range("Y16").Value = -0.68
$ws.Range("Z16").Value = 1.67
$ws.Range("W19").Value = 88.59
$ws.Range("X19").Value = 74.14
$ws.Range("Y19").Value = 91.8
$ws.Range("Z19").Value = 74.02
$ws.Range("W21").Value = -6.47
$ws.Range("X21").Value = -4.97
$ws.Range("Y21").Value = -5.97
$ws.Range("Z21").Value = -3.85
$ws.Range("W24").Value = 88.59
$ws.Range("X24").Value = 74.14
$ws.Range("Y24").Value = 91.8
$ws.Range("Z24").Value = 74.02
$ws.Range("W26").Value = -6.47
$ws.Range("X26").Value = -4.97
$ws.Range("Y26").Value = -5.97
$ws.Range("Z26").Value = -3.85
$ws.Range("W29").Value = 88.01000000000001
$ws.Range("X29").Value = 73.79000000000001
$ws.Range("Y29").Value = 91.29000000000001
$ws.Range("Z29").Value = 73.54000000000001
$ws.Range("W31").Value = -7.04
$ws.Range("X31").Value = -5.31
$ws.Range("Y31").Value = -6.48
$ws.Range("Z31").Value = -4.34
$ws.Range("W34").Value = 93.92
$ws.Range("X34").Value = 79.5
$ws.Range("Y34").Value = 100.48
$ws.Range("Z34").Value = 83.2
$ws.Range("W36").Value = -1.13
$ws.Range("X36").Value = 0.4
$ws.Range("Y36").Value = 2.71
$ws.Range("Z36").Value = 5.32
$ws.Range("W39").Value = 88.92
$ws.Range("X39").Value = 74.20999999999999
$ws.Range("Y39").Value = 91.89
$ws.Range("Z39").Value = 74.23999999999999
$ws.Range("W41").Value = -6.14
$ws.Range("X41").Value = -4.9
$ws.Range("Y41").Value = -5.88
$ws.Range("Z41").Value = -3.64
$ws.Range("W44").Value = 92.81999999999999
$ws.Range("X44").Value = 77.03
$ws.Range("Y44").Value = 96.04000000000001
$ws.Range("Z44").Value = 77.56
$ws.Range("W46").Value = -2.23
$ws.Range("X46").Value = -2.08
$ws.Range("Y46").Value = -1.73
$ws.Range("Z46").Value = -0.31
$ws.Range("W49").Value = 101.23
$ws.Range("X49").Value = 83.27
$ws.Range("Y49").Value = 102.59
$ws.Range("Z49").Value = 81.97
$ws.Range("W51").Value = 6.17
$ws.Range("X51").Value = 4.16
$ws.Range("Y51").Value = 4.82
$ws.Range("Z51").Value = 4.1
$ws.Range("W54").Value = 93.83
$ws.Range("X54").Value = 80.06999999999999
$ws.Range("Y54").Value = 99.45999999999999
$ws.Range("Z54").Value = 79.95
$ws.Range("W56").Value = -1.22
$ws.Range("X56").Value = 0.96
$ws.Range("Y56").Value = 1.69
$ws.Range("Z56").Value = 2.08
$ws.Range("W59").Value = 99.53
$ws.Range("X59").Value = 82.23
$ws.Range("Y59").Value = 103.24
$ws.Range("Z59").Value = 82.93000000000001
$ws.Range("W61").Value = 4.48
$ws.Range("X61").Value = 3.12
$ws.Range("Y61").Value = 5.47
$ws.Range("Z61").Value = 5.06
$ws.Range("W64").Value = 101.66
$ws.Range("X64").Value = 83.89
$ws.Range("Y64").Value = 105.24
$ws.Range("Z64").Value = 84.45999999999999
$ws.Range("W66").Value = 6.61
$ws.Range("X66").Value = 4.78
$ws.Range("Y66").Value = 7.47
$ws.Range("Z66").Value = 6.59
$ws.Range("W69").Value = 101.88
$ws.Range("X69").Value = 83.89
$ws.Range("Y69").Value = 106.27
$ws.Range("Z69").Value = 85.58
$ws.Range("W71").Value = 6.83
$ws.Range("X71").Value = 4.78
$ws.Range("Y71").Value = 8.5
$ws.Range("Z71").Value = 7.7
$ws.Range("W74").Value = 100.69
$ws.Range("X74").Value = 83.09
$ws.Range("Y74").Value = 104.12
$ws.Range("Z74").Value = 83.56
$ws.Range("W76").Value = 5.64
$ws.Range("X76").Value = 3.99
$ws.Range("Y76").Value = 6.35
$ws.Range("Z76").Value = 5.68
$ws.Range("W79").Value = 101.1
$ws.Range("X79").Value = 83.51000000000001
$ws.Range("Y79").Value = 104.57
$ws.Range("Z79").Value = 83.89
$ws.Range("W81").Value = 6.05
$ws.Range("X81").Value = 4.41
$ws.Range("Y81").Value = 6.8
$ws.Range("Z81").Value = 6.02
$ws.Range("W84").Value = 91.84
$ws.Range("X84").Value = 80.39
$ws.Range("Y84").Value = 99.77
$ws.Range("Z84").Value = 80.12
$ws.Range("W86").Value = -3.21
$ws.Range("X86").Value = 1.29
$ws.Range("Y86").Value = 2
$ws.Range("Z86").Value = 2.24
$ws.Range("W89").Value = 88.01000000000001
$ws.Range("X89").Value = 73.79000000000001
$ws.Range("Y89").Value = 91.29000000000001
$ws.Range("Z89").Value = 73.54000000000001
$ws.Range("W91").Value = -7.04
$ws.Range("X91").Value = -5.31
$ws.Range("Y91").Value = -6.48
$ws.Range("Z91").Value = -4.34
